$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "Status" column before the existing column C (vNameAccount),
# shifting C:G -> D:H. Excel copies the left neighbour's formatting into the
# freshly inserted column.
$ws.Columns("C").Insert()

# New "Status" column contents.
$ws.Range("C1").Value = "Status"
$ws.Range("C2").Value = "Failed"
$ws.Range("C3").Value = "Failed"
$ws.Range("C4").Value = "Failed"
$ws.Range("C5").Value = ""
$ws.Range("C6").Value = ""

# CT 02 / CT 03 rows now report "No" instead of "Yes".
$ws.Range("B3").Value = "No"
$ws.Range("B4").Value = "No"

# Refresh the "last run" timestamps (old column G, now column H) for the
# first three test cases; rows 5 and 6 keep their previous date stamp.
$ws.Range("H2").Value = "25_04_2020--19_25_34 264"
$ws.Range("H3").Value = "25_04_2020--19_14_17 408"
$ws.Range("H4").Value = "25_04_2020--19_14_48 020"

# Match the column widths captured in the saved workbook.
$ws.Columns("C").ColumnWidth = 9
$ws.Columns("H").ColumnWidth = 24.5

$null = $ws.Range("B5").Select()
